# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Updates the "K" column (column G) values on the active sheet for rows 2-20
# with the newly-calculated strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 7
    4  = 6
    5  = 4
    6  = 10
    7  = 8
    8  = 6
    9  = 2
    10 = 4
    11 = 11
    12 = 7
    13 = 7
    14 = 5
    15 = 4
    16 = 2
    17 = 5
    18 = 5
    19 = 4
    20 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
